$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 773, shifting existing rows 773-814 down to 774-815
$ws.Rows.Item(773).Insert()

# Populate the newly inserted row 773 with the new data.
# The date-like text in column A must stay a literal text string (matching
# the rest of the column), not get auto-converted into a date serial
# number. Writing it as a text formula and then pasting the computed
# value back over itself keeps the cell as plain text without picking up
# an extra number-format style.
$cellA = $ws.Cells.Item(773, 1)
$cellA.Formula = "=""2026/02/08"""
$cellA.Copy()
$cellA.PasteSpecial(-4163)

$ws.Cells.Item(773, 2).Value = "日"
$ws.Cells.Item(773, 3).Value = 17
$ws.Cells.Item(773, 4).Value = 94
